# ADT_A04 Channel with in built Mappers
#
# Renames the "A01_Patient_Delago" sheet to "A04_Patient_Delago" (and, as a
# consequence, the dependent "_xlnm._FilterDatabase" defined name updates
# automatically), then restores the view/selection state recorded for each
# sheet: the patient-message sheet ends up scrolled/selected at B28 and the
# "Table" lookup sheet ends up with D6:E6 selected.

$wb = $excel.ActiveWorkbook

# --- Rename the main sheet: A01_Patient_Delago -> A04_Patient_Delago -------
$mainSheet = $wb.Worksheets.Item("A01_Patient_Delago")
$mainSheet.Name = "A04_Patient_Delago"

# --- "Table" sheet: move the selection to D6:E6 ----------------------------
$tableSheet = $wb.Worksheets.Item("Table")
$tableSheet.Activate()
$tableSheet.Range("D6:E6").Select()

# --- Main sheet: scroll/select B28, keep it the active tab -----------------
$mainSheet.Activate()
$mainSheet.Range("B28").Select()
